$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

# Remove the "Life cycle cost" (row 21) and "Road safety impacts" (row 22) rows,
# which pushes all later rows up by two.
$ws.Rows.Item(21).Resize(2).Delete()

# The hidden "_FilterDatabase" name still points at the old (now too-large) range;
# shrink it to match the new data extent.
foreach ($n in $wb.Names) {
  if ($n.Name -eq "Data!_FilterDatabase") {
    $n.RefersTo = "=Data!`$A`$1:`$I`$23"
  }
}

# Restore the view: frozen pane/scroll position and active cell as left after
# the edit (near the bottom of the now-shorter table).
$ws.Activate() | Out-Null
$excel.ActiveWindow.ScrollRow = 14
$ws.Range("A22").Select() | Out-Null
